$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.090.54'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.81'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.6899'
$ws.Range("E5").Value = '  -6.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '236.72'
$ws.Range("E6").Value = '  -2.33%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3033'
$ws.Range("E8").Value = '  -3.65%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07527'
$ws.Range("E9").Value = '  +4.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.23'
$ws.Range("E10").Value = '  -5.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08071'
$ws.Range("E11").Value = '  -2.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.878.09'
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7201'
$ws.Range("E13").Value = '  -3.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.171'
$ws.Range("E14").Value = '  -3.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '88.49'
$ws.Range("E15").Value = '  -3.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.253.31'
$ws.Range("E16").Value = '  -2.12%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.763'
$ws.Range("E17").Value = '  -5.54%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '240.97'
$ws.Range("E18").Value = '  -2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007660'
$ws.Range("E19").Value = '  -2.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.99'
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.118.27'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.606'
$ws.Range("E24").Value = '  -4.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.980'
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.90'
$ws.Range("E26").Value = '  -2.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1457'
$ws.Range("E27").Value = '  -5.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.00'
$ws.Range("E28").Value = '  -3.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.922'
$ws.Range("E29").Value = '  -4.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.373'
$ws.Range("E30").Value = '  -8.15%  '
$ws.Range("E31").Value = '  -3.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.488'
$ws.Range("E32").Value = '  -3.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.037'
$ws.Range("E33").Value = '  -4.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05204'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.180'
$ws.Range("E35").Value = '  -4.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7101'
$ws.Range("E36").Value = '  -5.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9975'
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.661'
$ws.Range("E38").Value = '  -1.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01852'
$ws.Range("E39").Value = '  -5.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.679'
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9134'
$ws.Range("E41").Value = '  +5.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.932'
$ws.Range("E42").Value = '  -3.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4265'
$ws.Range("E43").Value = '  -5.53%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.050.33'
$ws.Range("E44").Value = '  -5.61%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.73'
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.25'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.008.38'
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.156'
$ws.Range("E49").Value = '  -6.01%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.737'
$ws.Range("E50").Value = '  -6.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.192'
$ws.Range("E51").Value = '  -3.39%  '
